$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1) Sprint Backlog sheet ("Sprint 2" follow-up tasks added to the table)
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Sprint Backlog")

# Set cell text values in the exact order the author originally typed them
# so new shared-string entries land in the same order as the target file.
$ws.Cells.Item(11, 3).Value = "Datenbankanbindung von der Medikation"
$ws.Cells.Item(14, 3).Value = "Datenbankanbindung von den Patienten"
$ws.Cells.Item(15, 3).Value = "Patientendaten genierieren"
$ws.Cells.Item(16, 3).Value = "View Anpassen nach Wunsch des Product Owner"
$ws.Cells.Item(17, 3).Value = "JPA-Anbindung"
$ws.Cells.Item(18, 3).Value = "MVP"
$ws.Cells.Item(18, 4).Value = "Alle Klassen nach MVP implementieren"
$ws.Cells.Item(13, 3).Value = "Datenbankanbindung von Reporten"
$ws.Cells.Item(12, 3).Value = "Datenbankanbindung von Description"

# Sprint number (column B) for every new row
2, 2, 2, 2, 2, 2, 2, 2 | Out-Null
for ($r = 11; $r -le 18; $r++) {
    $b = $ws.Cells.Item($r, 2)
    $b.Value = 2
    $b.HorizontalAlignment = -4108
    $b.VerticalAlignment = -4108
}

# Status (column L) = "waiting" for every new row
for ($r = 11; $r -le 18; $r++) {
    $l = $ws.Cells.Item($r, 12)
    $l.Value = "waiting"
    $l.HorizontalAlignment = -4108
    $l.VerticalAlignment = -4108
}

# Column C (Name) formatting per row
$c11 = $ws.Cells.Item(11, 3)
$c11.HorizontalAlignment = -4108
$c11.VerticalAlignment = -4108
$c11.WrapText = $true

foreach ($r in 12, 13, 14, 16) {
    $c = $ws.Cells.Item($r, 3)
    $c.HorizontalAlignment = -4108
    $c.WrapText = $true
}

foreach ($r in 15, 17, 18) {
    $c = $ws.Cells.Item($r, 3)
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
}

# Column D (Description) - wrap-formatted, mostly empty except row 18
foreach ($r in 12, 13, 14, 15, 16, 17, 18) {
    $ws.Cells.Item($r, 4).WrapText = $true
}

# Row heights
$ws.Rows.Item(11).RowHeight = 41
$ws.Rows.Item(12).RowHeight = 43
$ws.Rows.Item(13).RowHeight = 32
$ws.Rows.Item(14).RowHeight = 34
$ws.Rows.Item(15).RowHeight = 24
$ws.Rows.Item(16).RowHeight = 31
$ws.Rows.Item(17).RowHeight = 21
$ws.Rows.Item(18).RowHeight = 21

# Extend the "Tabelle1" table / autofilter range to cover the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:L18"))

# Sheet view: Sprint Backlog becomes the active tab, scrolled/selected per target
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("D12").Select()
$excel.ActiveWindow.Zoom = 138

# ----------------------------------------------------------------------
# 2) BurndownChart sheet: no longer the active tab, view scrolled up
# ----------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("BurndownChart")
$ws4.Range("E12").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.Zoom = 100

# Re-activate Sprint Backlog so it is the selected/visible tab on open
$ws.Activate()
